# "11 May - Noche"
#
# 1) The student previously listed as "SALAZAR MARIA FERNANDA ROSAS" had her
#    Paterno / Materno / Nombres mixed up. She is corrected to
#    "ROSAS SALAZAR MARIA FERNANDA" (Paterno=ROSAS, Materno=SALAZAR,
#    Nombres=MARIA FERNANDA).
#
# 2) Because "Calificaciones" is kept sorted alphabetically by full name, and
#    "ROSAS ..." now sorts before "RUIZ LOPEZ XIMENA MICHELL" (previously the
#    row right above her), the two rows swap places: row 36 (was RUIZ) and
#    row 37 (was SALAZAR...) exchange their full contents (name + all grade
#    columns).
#
# 3) In "Totales Blanco" the same student's row keeps its row position (that
#    sheet isn't name-sorted) but the Paterno/Materno/Nombres cells are
#    rotated to their corrected values.

$wb = $excel.ActiveWorkbook

# --- 1) & 2): swap rows 36 / 37 on "Calificaciones" -------------------------
$wsCal = $wb.Worksheets.Item("Calificaciones")

$row36 = $wsCal.Range("A36:S36").Value2
$row37 = $wsCal.Range("A37:S37").Value2

$wsCal.Range("A36:S36").Value2 = $row37
$wsCal.Range("A37:S37").Value2 = $row36

# Name column keeps its corrected spelling explicitly (belt & suspenders on
# top of the row swap above, in case Value2 round-tripping ever changes the
# text of the copied array).
$wsCal.Cells.Item(36, 1).Value = "ROSAS SALAZAR MARIA FERNANDA"
$wsCal.Cells.Item(37, 1).Value = "RUIZ LOPEZ XIMENA MICHELL"

# --- 3): fix Paterno/Materno/Nombres rotation on "Totales Blanco" ----------
$wsBlanco = $wb.Worksheets.Item("Totales Blanco")

$wsBlanco.Cells.Item(37, 2).Value = "ROSAS"
$wsBlanco.Cells.Item(37, 3).Value = "SALAZAR"
$wsBlanco.Cells.Item(37, 4).Value = "MARIA FERNANDA"
